$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B:G
# Row 2 (Q0 -> label A2 stays "6")
$ws.Range("B2").Value = 0.1282474195497992
$ws.Range("C2").Value = 0.9763931104222339
$ws.Range("D2").Value = 4.446347975453152
$ws.Range("E2").Value = 2.108636520468417
$ws.Range("F2").Value = 2.125676030426229
$ws.Range("G2").Value = 51

# Row 3
$ws.Range("B3").Value = 0.1215397788006174
$ws.Range("C3").Value = 1.018530036818849
$ws.Range("D3").Value = 4.443793105881769
$ws.Range("E3").Value = 2.108030622614807
$ws.Range("F3").Value = 2.125890254825519
$ws.Range("G3").Value = 50

# Row 4
$ws.Range("B4").Value = 0.1466579445676272
$ws.Range("C4").Value = 0.9414795404354654
$ws.Range("D4").Value = 4.22668197375337
$ws.Range("E4").Value = 2.055889582091745
$ws.Range("F4").Value = 2.071902797105537
$ws.Range("G4").Value = 49

# Row 5
$ws.Range("B5").Value = 0.1692583036833795
$ws.Range("C5").Value = 1.044790455529975
$ws.Range("D5").Value = 4.665004838604097
$ws.Range("E5").Value = 2.159862226764498
$ws.Range("F5").Value = 2.176006054099445
$ws.Range("G5").Value = 48

# Row 6
$ws.Range("B6").Value = 0.1384867545238795
$ws.Range("C6").Value = 0.9904347050370882
$ws.Range("D6").Value = 4.525272260743743
$ws.Range("E6").Value = 2.127268732610843
$ws.Range("F6").Value = 2.145705533808543
$ws.Range("G6").Value = 47

# Row 7
$ws.Range("B7").Value = 0.1625861655212504
$ws.Range("C7").Value = 1.003062538999909
$ws.Range("D7").Value = 4.667669780972878
$ws.Range("E7").Value = 2.160479062840665
$ws.Range("F7").Value = 2.17815841638289
$ws.Range("G7").Value = 46

# Row 8
$ws.Range("B8").Value = 0.09841386822182357
$ws.Range("C8").Value = 0.9352685283506139
$ws.Range("D8").Value = 4.357710809747646
$ws.Range("E8").Value = 2.087513068162124
$ws.Range("F8").Value = 2.108754201445493
$ws.Range("G8").Value = 45

# Row 9
$ws.Range("B9").Value = 0.06748005575673845
$ws.Range("C9").Value = 0.9265491817941878
$ws.Range("D9").Value = 4.470640045755423
$ws.Range("E9").Value = 2.114388811395724
$ws.Range("F9").Value = 2.137743895861163
$ws.Range("G9").Value = 44

# Row 10
$ws.Range("B10").Value = 0.09914051671729636
$ws.Range("C10").Value = 0.9226436494409858
$ws.Range("D10").Value = 4.481900347334963
$ws.Range("E10").Value = 2.117049916117937
$ws.Range("F10").Value = 2.13975451822466
$ws.Range("G10").Value = 43

# Row 11
$ws.Range("B11").Value = 0.07576832823998672
$ws.Range("C11").Value = 0.9176283451160873
$ws.Range("D11").Value = 4.557762439856401
$ws.Range("E11").Value = 2.13489166934915
$ws.Range("F11").Value = 2.159408835161238
$ws.Range("G11").Value = 42

$wb.Save()
